$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2 through 396)
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07).
for ($row = 2; $row -le 396; $row++) {
    $ws.Cells.Item($row, 3).Value = 45206
}
